# Update cryptos list (price + 1h volume change) columns D and E
# per the GitHub Actions data refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Cell, $Text)
    # Force text storage (keep formats like "1.00" or "67.124.49"
    # from being reinterpreted as numbers), then restore the
    # default "Normal" style so no stray formatting is introduced.
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "67.124.49"
Set-TextCell $ws.Range("E2") "  +4.59%  "
Set-TextCell $ws.Range("D3") "3.464.91"
Set-TextCell $ws.Range("E3") "  +4.25%  "
Set-TextCell $ws.Range("D4") "1.00"
Set-TextCell $ws.Range("E4") "  +0.00%  "
Set-TextCell $ws.Range("D5") "587.55"
Set-TextCell $ws.Range("E5") "  +6.48%  "
Set-TextCell $ws.Range("D6") "187.38"
Set-TextCell $ws.Range("E6") "  +8.61%  "
Set-TextCell $ws.Range("E7") "  +1.13%  "
Set-TextCell $ws.Range("D8") "3.458.62"
Set-TextCell $ws.Range("E8") "  +4.37%  "
Set-TextCell $ws.Range("E9") "  -0.06%  "
Set-TextCell $ws.Range("E10") "  +0.25%  "
Set-TextCell $ws.Range("E11") "  +2.53%  "
Set-TextCell $ws.Range("D12") "56.06"
Set-TextCell $ws.Range("E12") "  +5.71%  "
Set-TextCell $ws.Range("D13") "0.0000279"
Set-TextCell $ws.Range("E13") "  +0.71%  "
Set-TextCell $ws.Range("D14") "9.41"
Set-TextCell $ws.Range("E14") "  +4.15%  "
Set-TextCell $ws.Range("D15") "4.021.07"
Set-TextCell $ws.Range("E15") "  +4.47%  "
Set-TextCell $ws.Range("E16") "  +4.15%  "
Set-TextCell $ws.Range("D17") "3.460.24"
Set-TextCell $ws.Range("E17") "  +4.49%  "
Set-TextCell $ws.Range("D18") "67.088.24"
Set-TextCell $ws.Range("E18") "  +4.69%  "
Set-TextCell $ws.Range("D19") "12.15"
Set-TextCell $ws.Range("E19") "  +4.26%  "
Set-TextCell $ws.Range("D20") "0.118"
Set-TextCell $ws.Range("E20") "  -1.53%  "
Set-TextCell $ws.Range("E21") "  +3.77%  "
Set-TextCell $ws.Range("D22") "485.67"
Set-TextCell $ws.Range("E22") "  +7.78%  "
Set-TextCell $ws.Range("D23") "5.28"
Set-TextCell $ws.Range("E23") "  +5.78%  "
Set-TextCell $ws.Range("D24") "16.86"
Set-TextCell $ws.Range("E24") "  +21.16%  "
Set-TextCell $ws.Range("E25") "  +11.48%  "
Set-TextCell $ws.Range("D26") "89.57"
Set-TextCell $ws.Range("E26") "  +3.00%  "
Set-TextCell $ws.Range("D27") "2.95"
Set-TextCell $ws.Range("E27") "  +3.41%  "
Set-TextCell $ws.Range("D28") "10.95"
Set-TextCell $ws.Range("E28") "  +3.91%  "
Set-TextCell $ws.Range("D29") "9.10"
Set-TextCell $ws.Range("E29") "  +6.35%  "
Set-TextCell $ws.Range("D30") "31.36"
Set-TextCell $ws.Range("E30") "  +1.90%  "
Set-TextCell $ws.Range("D31") "7.17"
Set-TextCell $ws.Range("E31") "  +10.81%  "
Set-TextCell $ws.Range("D32") "598.58"
Set-TextCell $ws.Range("E32") "  +5.19%  "
Set-TextCell $ws.Range("E33") "  +3.36%  "
Set-TextCell $ws.Range("D34") "63.98"
Set-TextCell $ws.Range("E34") "  +2.16%  "
Set-TextCell $ws.Range("D35") "0.112"
Set-TextCell $ws.Range("E35") "  +5.34%  "
Set-TextCell $ws.Range("D36") "0.150"
Set-TextCell $ws.Range("E36") "  +6.03%  "
Set-TextCell $ws.Range("D38") "36.53"
Set-TextCell $ws.Range("E38") "  +4.06%  "
Set-TextCell $ws.Range("D39") "3.55"
Set-TextCell $ws.Range("E39") "  +1.06%  "
Set-TextCell $ws.Range("E40") "  +5.51%  "
Set-TextCell $ws.Range("E41") "  +4.49%  "
Set-TextCell $ws.Range("D42") "3.239.63"
Set-TextCell $ws.Range("E42") "  +5.98%  "
Set-TextCell $ws.Range("E43") "  +7.32%  "
Set-TextCell $ws.Range("E44") "  +4.60%  "
Set-TextCell $ws.Range("E45") "  +4.56%  "
Set-TextCell $ws.Range("E46") "  +3.65%  "
Set-TextCell $ws.Range("D47") "2.76"
Set-TextCell $ws.Range("E47") "  +23.86%  "
Set-TextCell $ws.Range("D48") "0.135"
Set-TextCell $ws.Range("E48") "  +2.05%  "
Set-TextCell $ws.Range("D49") "3.28"
Set-TextCell $ws.Range("E49") "  +14.29%  "
Set-TextCell $ws.Range("D50") "8.75"
Set-TextCell $ws.Range("E50") "  +7.48%  "
Set-TextCell $ws.Range("D51") "1.00"
Set-TextCell $ws.Range("E51") "  +0.09%  "
